# "update allowing dynamic number of stages. not tested"
#
# The author inserted a new (blank) row at the top of the "cost1" sheet -
# shifting its existing content down by one row - to make room for
# additional rows/stages going forward, then left the workbook with the
# "cost1" sheet active and cell B11 selected (mid-edit, matching the
# commit message "not tested").

$wb = $excel.ActiveWorkbook

# --- cost1: insert a blank row above row 1 (data shifts down to rows 2-7) ---
$costSheet = $wb.Worksheets.Item("cost1")
$costSheet.Rows.Item(1).Insert()

# --- defined names: Solver's ranges point at a sheet ("Elmore's Ski
#     Boots") that lives in another workbook, not one of the local sheets;
#     mark them as external (workbook index [1]) references ---
foreach ($n in $wb.Names) {
    if ($n.RefersTo -like "*'Elmore''s Ski Boots'*" -and $n.RefersTo -notlike "*'[1]Elmore''s Ski Boots'*") {
        $n.RefersTo = $n.RefersTo -replace "'Elmore''s Ski Boots'", "'[1]Elmore''s Ski Boots'"
    }
}

# --- leave the workbook focused on "cost1" with B11 selected/zoomed in,
#     and drop the old "Main" tab selection/zoom it had before ---
$mainSheet = $wb.Worksheets.Item("Main")
$mainSheet.Activate()
$excel.ActiveWindow.Zoom = 208

$costSheet.Activate()
$excel.ActiveWindow.Zoom = 255
$costSheet.Range("B11").Select()
